$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 484 (after row 483), shifting the rest of the
# data (old rows 484-497) down to rows 487-500.
$ws.Rows.Item(484).Resize(3).Insert()

# New row 484: new weekly entry (Primera) for date 45041
$ws.Cells.Item(484, 1).Value = 4
$ws.Cells.Item(484, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(484, 3).Value = "Los Lagos"
$ws.Cells.Item(484, 4).Value = 45041
$ws.Cells.Item(484, 5).Value = 10
$ws.Cells.Item(484, 6).Value = "Fruta"
$ws.Cells.Item(484, 7).Value = 100101
$ws.Cells.Item(484, 8).Value = "Berries"
$ws.Cells.Item(484, 9).Value = 100101007
$ws.Cells.Item(484, 10).Value = "Kiwi"
$ws.Cells.Item(484, 11).Value = "Hayward"
$ws.Cells.Item(484, 12).Value = "Primera"
$ws.Cells.Item(484, 13).Value = 300
$ws.Cells.Item(484, 14).Value = 18000
$ws.Cells.Item(484, 15).Value = 18000
$ws.Cells.Item(484, 16).Value = 18000
$ws.Cells.Item(484, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(484, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(484, 19).Value = 1200
$ws.Cells.Item(484, 20).Value = 15

# New row 485: new weekly entry (Segunda) for date 45041
$ws.Cells.Item(485, 1).Value = 4
$ws.Cells.Item(485, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(485, 3).Value = "Los Lagos"
$ws.Cells.Item(485, 4).Value = 45041
$ws.Cells.Item(485, 5).Value = 10
$ws.Cells.Item(485, 6).Value = "Fruta"
$ws.Cells.Item(485, 7).Value = 100101
$ws.Cells.Item(485, 8).Value = "Berries"
$ws.Cells.Item(485, 9).Value = 100101007
$ws.Cells.Item(485, 10).Value = "Kiwi"
$ws.Cells.Item(485, 11).Value = "Hayward"
$ws.Cells.Item(485, 12).Value = "Segunda"
$ws.Cells.Item(485, 13).Value = 300
$ws.Cells.Item(485, 14).Value = 15000
$ws.Cells.Item(485, 15).Value = 15000
$ws.Cells.Item(485, 16).Value = 15000
$ws.Cells.Item(485, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(485, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(485, 19).Value = 1000
$ws.Cells.Item(485, 20).Value = 15

# New row 486: duplicate of (original/current) row 483 (Especial, date 44705)
$ws.Cells.Item(486, 1).Value = 4
$ws.Cells.Item(486, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(486, 3).Value = "Los Lagos"
$ws.Cells.Item(486, 4).Value = 44705
$ws.Cells.Item(486, 5).Value = 10
$ws.Cells.Item(486, 6).Value = "Fruta"
$ws.Cells.Item(486, 7).Value = 100101
$ws.Cells.Item(486, 8).Value = "Berries"
$ws.Cells.Item(486, 9).Value = 100101007
$ws.Cells.Item(486, 10).Value = "Kiwi"
$ws.Cells.Item(486, 11).Value = "Hayward"
$ws.Cells.Item(486, 12).Value = "Especial"
$ws.Cells.Item(486, 13).Value = 300
$ws.Cells.Item(486, 14).Value = 20000
$ws.Cells.Item(486, 15).Value = 20000
$ws.Cells.Item(486, 16).Value = 20000
$ws.Cells.Item(486, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(486, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(486, 19).Value = 1333
$ws.Cells.Item(486, 20).Value = 15

# Copy the date formatting/style from row 483 column D to the new rows' column D
$ws.Cells.Item(483, 4).Copy()
$ws.Range($ws.Cells.Item(484, 4), $ws.Cells.Item(486, 4)).PasteSpecial(-4122)
$excel.CutCopyMode = 0
